$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E hold numeric/percentage-looking text that must remain text.
# Force Text number format on the full D:E data range before assigning values so
# Excel does not auto-convert strings like "1.006" or "29.146.75" into numbers,
# then reset the style back to Normal (default) afterwards so cell styling is
# unchanged from the original (no explicit "s" style index).
$ws.Range("D2:E51").NumberFormat = "@"

# Coin name / link swaps (rows 46-48 re-ranked)
$ws.Range('B46').Value = 'XinFinNetwork'
$ws.Range('C46').Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('B48').Value = 'TheSandbox'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'

# Updated prices and volume percentages
$ws.Range('D2').Value = '29.146.75'
$ws.Range('E2').Value = '  +0.43%  '
$ws.Range('D3').Value = '1.836.13'
$ws.Range('E3').Value = '  +0.20%  '
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  +0.71%  '
$ws.Range('D5').Value = '244.12'
$ws.Range('E5').Value = '  +1.12%  '
$ws.Range('D6').Value = '0.6278'
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('D7').Value = '1.006'
$ws.Range('E7').Value = '  +0.65%  '
$ws.Range('D8').Value = '0.07462'
$ws.Range('E8').Value = '  -1.82%  '
$ws.Range('D9').Value = '0.2930'
$ws.Range('E9').Value = '  +0.52%  '
$ws.Range('D10').Value = '23.01'
$ws.Range('E10').Value = '  +1.31%  '
$ws.Range('D11').Value = '0.07733'
$ws.Range('E11').Value = '  +0.00%  '
$ws.Range('D12').Value = '1.833.49'
$ws.Range('E12').Value = '  -0.05%  '
$ws.Range('D13').Value = '5.004'
$ws.Range('E13').Value = '  +1.00%  '
$ws.Range('D14').Value = '0.6672'
$ws.Range('E14').Value = '  +0.53%  '
$ws.Range('D15').Value = '83.07'
$ws.Range('E15').Value = '  +0.56%  '
$ws.Range('D16').Value = '0.000009341'
$ws.Range('E16').Value = '  -4.30%  '
$ws.Range('D17').Value = '6.056'
$ws.Range('E17').Value = '  +1.15%  '
$ws.Range('D18').Value = '29.159.16'
$ws.Range('E18').Value = '  +0.49%  '
$ws.Range('E19').Value = '  +2.23%  '
$ws.Range('D20').Value = '223.58'
$ws.Range('E20').Value = '  -1.01%  '
$ws.Range('D21').Value = '1.007'
$ws.Range('E21').Value = '  +0.77%  '
$ws.Range('D22').Value = '7.121'
$ws.Range('E22').Value = '  -1.26%  '
$ws.Range('D23').Value = '1.007'
$ws.Range('E23').Value = '  +0.74%  '
$ws.Range('D24').Value = '160.30'
$ws.Range('E24').Value = '  +1.36%  '
$ws.Range('D25').Value = '0.1403'
$ws.Range('E25').Value = '  +2.25%  '
$ws.Range('D26').Value = '8.500'
$ws.Range('E26').Value = '  +1.03%  '
$ws.Range('D27').Value = '17.92'
$ws.Range('E27').Value = '  +0.34%  '
$ws.Range('D28').Value = '1.501'
$ws.Range('D29').Value = '4.147'
$ws.Range('E29').Value = '  +2.23%  '
$ws.Range('D30').Value = '4.064'
$ws.Range('E30').Value = '  +1.05%  '
$ws.Range('D31').Value = '0.05464'
$ws.Range('E31').Value = '  +5.38%  '
$ws.Range('E32').Value = '  +0.62%  '
$ws.Range('D33').Value = '0.7483'
$ws.Range('E33').Value = '  +1.31%  '
$ws.Range('D34').Value = '1.851'
$ws.Range('E34').Value = '  +0.24%  '
$ws.Range('E35').Value = '  -0.76%  '
$ws.Range('D36').Value = '2.616'
$ws.Range('E36').Value = '  -2.95%  '
$ws.Range('D37').Value = '1.225.71'
$ws.Range('E37').Value = '  -3.21%  '
$ws.Range('D38').Value = '2.764'
$ws.Range('E38').Value = '  +0.07%  '
$ws.Range('D39').Value = '0.01787'
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('D40').Value = '6.652'
$ws.Range('E40').Value = '  +6.57%  '
$ws.Range('D41').Value = '0.8960'
$ws.Range('E41').Value = '  +0.17%  '
$ws.Range('D42').Value = '1.006'
$ws.Range('E42').Value = '  +0.66%  '
$ws.Range('D43').Value = '101.56'
$ws.Range('E43').Value = '  +0.22%  '
$ws.Range('D44').Value = '65.40'
$ws.Range('E44').Value = '  +1.49%  '
$ws.Range('E45').Value = '  -1.96%  '
$ws.Range('D46').Value = '0.07755'
$ws.Range('E46').Value = '  +13.85%  '
$ws.Range('D47').Value = '0.5114'
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('D48').Value = '0.4042'
$ws.Range('E48').Value = '  +1.52%  '
$ws.Range('D49').Value = '8.949'
$ws.Range('E49').Value = '  +1.01%  '
$ws.Range('D50').Value = '0.05814'
$ws.Range('E50').Value = '  +1.10%  '
$ws.Range('D51').Value = '1.648'
$ws.Range('E51').Value = '  +1.62%  '

# Restore default styling on the D:E range (clears the temporary Text format)
$ws.Range("D2:E51").Style = "Normal"

